$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Content.Find.Execute("2024-01-09 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-01-10 Wednesday", 2)

# Update the division-problem table cells by direct (row, column) address so
# that no partially-replaced text can be re-matched by a later lookup
# (several of the new values coincide with other cells' original values).
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("50÷7=", "24÷7=", "66÷5=", "17÷5=", "82÷4=")
    5  = @("33÷8=", "15÷9=", "78÷2=", "12÷3=", "64÷2=")
    9  = @("60÷7=", "56÷2=", "19÷2=", "91÷5=", "62÷6=")
    13 = @("62÷5=", "71÷2=", "89÷7=", "99÷5=", "97÷6=")
    17 = @("41÷5=", "32÷3=", "39÷5=", "33÷8=", "91÷9=")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le $values.Count; $col++) {
        $t.Cell($row, $col).Range.Text = $values[$col - 1]
    }
}

Write-Output "done"
